$d = $word.ActiveDocument

# Pull the full package XML (includes every part, not just document.xml)
$xml = $d.Content.XML()

# The portion of document.xml we need to replace: the bookmarked Heading1
# title paragraph followed by the bold "By Dorothy Day" paragraph.
$oldFragment = '<w:bookmarkStart w:id="0" w:name="suicide-or-sacrifice" /><w:p w14:paraId="00000001" w14:textId="77777777" w:rsidR="00000000" w:rsidRDefault="00000000"><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>Suicide or Sacrifice?</w:t></w:r></w:p><w:bookmarkEnd w:id="0" /><w:p w14:paraId="00000002" w14:textId="77777777" w:rsidR="00000000" w:rsidRDefault="00000000"><w:r><w:rPr><w:b/></w:rPr><w:t>By Dorothy Day</w:t></w:r></w:p>'

if (-not $xml.Contains($oldFragment)) {
    throw "old fragment not found"
}

$newFragment = '<w:p w14:paraId="00000001" w14:textId="77777777" w:rsidR="00000000" w:rsidRDefault="00000000"><w:pPr><w:pStyle w:val="Title"/></w:pPr><w:r><w:t xml:space="preserve">Suicide</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">or</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Sacrifice</w:t></w:r><w:r><w:t xml:space="preserve">?</w:t></w:r></w:p><w:p w14:paraId="00000002" w14:textId="77777777" w:rsidR="00000000" w:rsidRDefault="00000000"><w:pPr><w:pStyle w:val="Authors"/></w:pPr><w:r><w:t xml:space="preserve">Dorothy</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Day</w:t></w:r></w:p>'

$newXml = $xml.Replace($oldFragment, $newFragment)

$d.Content.InsertXML($newXml)
